$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Week #2's hours (Tue-Fri first, Monday/C6 last so the
# IF(C6="", 0, SUM(C6:I6)) totals formula recalculates against the
# already-written neighbour cells).
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1.5
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 2
$ws.Range("C6").Value = 0

$ws.Range("F6").Select()
